# Add new rows to the "Все персонажи/Торкоал" sheet: nine more
# SCRIPT/P02P01A/um####.ssb filename rows appended below the existing
# data (rows 7-15), matching the formatting already used by rows 4-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "SCRIPT/P02P01A/um1204.ssb",
    "SCRIPT/P02P01A/um1303.ssb",
    "SCRIPT/P02P01A/um1306.ssb",
    "SCRIPT/P02P01A/um1309.ssb",
    "SCRIPT/P02P01A/um1401.ssb",
    "SCRIPT/P02P01A/um1404.ssb",
    "SCRIPT/P02P01A/um1501.ssb",
    "SCRIPT/P02P01A/um1601.ssb",
    "SCRIPT/P02P01A/um1604.ssb"
)

$row = 7
foreach ($val in $newValues) {
    $cell = $ws.Range("A$row")
    $cell.Value = $val
    $cell.WrapText = $true
    $ws.Rows($row).RowHeight = 43.2
    $row++
}

# Scroll the view down and move the selection, matching where the
# author was working after adding the new rows.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C11").Select()
